# Generate Report for Handback
# Re-labels the two handback rows with freshly-generated UUID file names and
# timestamps, and consolidates both rows' generated .xlf artifact names/dates
# onto the single newest generation (row 2's), matching a fresh run of the
# handback-status report.

function Set-HyperlinkDisplay($ws, $addr, $newText) {
    $links = @($ws.Hyperlinks)
    foreach ($hl in $links) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $newText
        }
    }
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# New identifiers generated by this handback run
# ---------------------------------------------------------------------
$uuid1 = "3a299245-b98d-43ab-97fe-cfa3c405a008"
$uuid2 = "ffff1d7f2734-ff74-4121-aca5-0eaae420e7b2"
$hash  = "1a782bb7fb46b082a97797d5082ffc6a5365d026"

$md1 = "$uuid1.md"
$md2 = "$uuid2.md"

$xlfZhCn = "$uuid1.$hash.zh-cn.xlf"
$xlfDeDe = "$uuid1.$hash.de-de.xlf"

$genDate   = "2016-08-22 21:03:07"
$zhcnHoDate = "2016-08-22 21:02:57"
$zhcnHbDate = "2016-08-22 21:03:28"
$dedeHbDate = "2016-08-22 21:03:35"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = $md1
$wsOverview.Range("B2").Value = "e2e\$md1"
$wsOverview.Range("G2").Value = $genDate

$wsOverview.Range("A3").Value = $md2
$wsOverview.Range("B3").Value = "e2e\$md2"
$wsOverview.Range("G3").Value = $genDate

Set-HyperlinkDisplay $wsOverview "`$B`$2" "e2e\$md1"
Set-HyperlinkDisplay $wsOverview "`$B`$3" "e2e\$md2"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = $md1
$wsZhCn.Range("I2").Value = $md1
$wsZhCn.Range("A3").Value = $md2
$wsZhCn.Range("I3").Value = $md2

$wsZhCn.Range("G2").Value = $xlfZhCn
$wsZhCn.Range("J2").Value = $xlfZhCn
$wsZhCn.Range("H2").Value = $zhcnHoDate
$wsZhCn.Range("K2").Value = $zhcnHbDate

$wsZhCn.Range("G3").Value = $xlfZhCn
$wsZhCn.Range("J3").Value = $xlfZhCn
$wsZhCn.Range("H3").Value = $zhcnHoDate
$wsZhCn.Range("K3").Value = $zhcnHbDate

Set-HyperlinkDisplay $wsZhCn "`$A`$2" $md1
Set-HyperlinkDisplay $wsZhCn "`$I`$2" $md1
Set-HyperlinkDisplay $wsZhCn "`$A`$3" $md2
Set-HyperlinkDisplay $wsZhCn "`$I`$3" $md2

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = $md1
$wsDeDe.Range("I2").Value = $md1
$wsDeDe.Range("A3").Value = $md2
$wsDeDe.Range("I3").Value = $md2

$wsDeDe.Range("G2").Value = $xlfDeDe
$wsDeDe.Range("J2").Value = $xlfDeDe
$wsDeDe.Range("H2").Value = $genDate
$wsDeDe.Range("K2").Value = $dedeHbDate

$wsDeDe.Range("G3").Value = $xlfDeDe
$wsDeDe.Range("J3").Value = $xlfDeDe
$wsDeDe.Range("H3").Value = $genDate
$wsDeDe.Range("K3").Value = $dedeHbDate

Set-HyperlinkDisplay $wsDeDe "`$A`$2" $md1
Set-HyperlinkDisplay $wsDeDe "`$I`$2" $md1
Set-HyperlinkDisplay $wsDeDe "`$A`$3" $md2
Set-HyperlinkDisplay $wsDeDe "`$I`$3" $md2
